# Version 2 of Diplom
# Applies:
#  - teams: replace row3 (group/full_name/link_index) and drop rows 4-8
#  - lab1: renumber link_index values (A2: 1->5, A3: 2->1)
#  - lab4: add a new row (index=5, link=pipeline repo url)
#  - make "teams" the active sheet/tab again (was "lab5")

$wb = $excel.ActiveWorkbook

# ---- teams (sheet 1) ----
$teams = $wb.Worksheets.Item("teams")
$teams.Range("A3").Value = "РИМ-220909"
$teams.Range("B3").Value = "Иванов Иван Иванович"
$teams.Range("C3").Value = 2
$teams.Rows("4:8").Delete()

# ---- lab1 (sheet 2) ----
$lab1 = $wb.Worksheets.Item("lab1")
$lab1.Range("A2").Value = 5
$lab1.Range("A3").Value = 1

# ---- lab4 (sheet 5) ----
$lab4 = $wb.Worksheets.Item("lab4")
$lab4.Range("A3").Value = 5
$lab4.Range("B3").Value = "https://github.com/Wheatly99/The-simplest-machine-learning-pipeline.git"

# ---- restore selections / active sheet to match the authored state ----
$null = $lab1.Range("B14").Select()
$lab2 = $wb.Worksheets.Item("lab2")
$null = $lab2.Range("B2").Select()
$null = $lab4.Range("G9").Select()
$lab5 = $wb.Worksheets.Item("lab5")
$null = $lab5.Range("E14").Select()

$null = $teams.Range("D2:M4").Select()
$null = $teams.Activate()
